$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update DAMSLTag (column I) and DialogAct (column J) values for the rows
# re-annotated by re-running SGNN after transcript clean-up.
$updates = @(
    @{ Row = 22; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 23; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 34; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 35; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 42; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 48; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 53; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 70; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 104; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 122; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 124; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 144; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 145; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 182; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 184; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 186; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 209; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 212; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 213; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 220; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 224; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 226; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 239; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 240; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 242; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 246; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 255; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 257; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 265; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 271; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 274; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 281; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 287; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 290; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 291; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 303; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 312; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 314; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 322; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 342; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 344; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 346; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 349; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 356; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 359; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 361; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 364; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 369; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 381; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 394; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 395; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 398; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 401; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 418; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 432; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 433; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 437; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 444; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 446; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 456; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 462; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 475; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

Write-Output "Updated $($updates.Count) rows"
